$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91808
$ws.Range("B3").Value = 80377
$ws.Range("B4").Value = 83089
$ws.Range("B5").Value = 57884
